$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two "smart folder" claim/annulment numbers (numbers are entered
# with a leading apostrophe so Excel keeps them as text, preserving the
# leading zero and trailing spaces exactly like the original cells).
$ws.Range("F2").Value = "'0420194406833 "
$ws.Range("F3").Value = "'0420172008637  "

# Leave the active selection on F4, matching the saved cursor position
$ws.Range("F4").Select()
